$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Row 3 (header row) - only the D3 header text actually changes
# ("calculations" -> "values for calculation"); A3/B3/C3 stay the same.
# Written last among the "new text" cells so the shared-string table
# lands in the same append order Excel produced.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Row 4 - new first test case: "an unimplemented function"
# ---------------------------------------------------------------------
$ws.Range("C4").Value = "UnimpementedFunction Function"
$ws.Range("A4").Value = "an unimplemented function"
$ws.Range("D4").Value = 1
$ws.Range("B4").FormulaArray = "=TEXTSPLIT(D4,"","")"

# ---------------------------------------------------------------------
# Row 5 - the original vlooku/#NAME? test, shifted down from row 4
# ---------------------------------------------------------------------
$ws.Range("A5").Value = """#NAME?"""
$ws.Range("C5").Value = "FormulaParseException"
$ws.Range("B5").FormulaArray = "=vlooku"

# ---------------------------------------------------------------------
# Row 6 - WorkbookNotFoundException / #REF! (external reference formula)
# ---------------------------------------------------------------------
$ws.Range("A6").Value = """#REF!"""
$ws.Range("C6").Value = "WorkbookNotFoundException"
# First reference by file name so Excel registers the external link part...
$ws.Range("B6").Formula = "='[externalLink1.xlsx]Summary'!`$A`$1"
# ...then re-express it with the canonical [1] index form used by OOXML.
$ws.Range("B6").Formula = "=[1]Summary!`$A`$1"

# ---------------------------------------------------------------------
# Row 10 - Other error / ClassCastException (written here so the shared
# strings land in the same order the recorded workbook used)
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "Other error"
$ws.Range("C10").Value = "ClassCastException"

# ---------------------------------------------------------------------
# Row 7 - ClassCastException / #VALUE!
# ---------------------------------------------------------------------
$ws.Range("D7").Value = "あ"
$ws.Range("A7").Value = """#VALUE!"""
$ws.Range("C7").Value = "(No Exception)"
$ws.Range("B7").Formula = "=D7/2"

# ---------------------------------------------------------------------
# Row 8 - #DIV/0!
# ---------------------------------------------------------------------
$ws.Range("A8").Value = """#DIV/0!"""
$ws.Range("C8").Value = "(No Exception)"
$ws.Range("D8").Value = 0
$ws.Range("B8").Formula = "=1/D8"

# ---------------------------------------------------------------------
# Row 9 - #N/A (VLOOKUP miss)
# ---------------------------------------------------------------------
$ws.Range("A9").Value = """#N/A"""
$ws.Range("C9").Value = "(No Exception)"
$ws.Range("B9").Formula = "=VLOOKUP(""abc"",H3:I6,2,FALSE)"

# ---------------------------------------------------------------------
# Row 3 D3 header text update
# ---------------------------------------------------------------------
$ws.Range("D3").Value = "values for calculation"

# ---------------------------------------------------------------------
# Row 10 formulas (B10 array formula + D10 helper formula)
# ---------------------------------------------------------------------
$ws.Range("D10").Formula = "=VLOOKUP(""abc"",H3:I6,2,FALSE)"
$ws.Range("B10").FormulaArray = "=IFS(D10,""A"",TRUE,""B"")"

# ---------------------------------------------------------------------
# Column widths (best-fit, mirroring the widened A/C columns)
# ---------------------------------------------------------------------
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(3).AutoFit()

# ---------------------------------------------------------------------
# Selection moves to A11
# ---------------------------------------------------------------------
$ws.Range("A11").Select()
